# Update calibration estimates with std moving average values.
$wb = $excel.ActiveWorkbook

# --- Sheet "Linear" ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.0003262946206169515
$wsLinear.Range("B3").Value = -0.002694509524860991
$wsLinear.Range("B4").Value = 0.02787027567024508
$wsLinear.Range("B5").Value = "[1.0, 0.2250619841574918, 0.05663358377404991, 0.08601009109933162, 0.04700098638789774, 0.055161362408366214, 0.20023958932650116, 0.3616852570962582, 0.21141811756105128, 0.04407705268496192, 0.03818347873613333, 0.04385984419680644, 0.06355684765771173, 0.20225998594406605, 0.36498165762185153, 0.205863827717201, 0.019719607061920656, 0.04427262455289715, 0.03756950390983698, 0.037981247201255335]"

# --- Sheet "NonLinear" ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.003190840108231945
$wsNonLinear.Range("B5").Value = 0.002846808004166386
$wsNonLinear.Range("B6").Value = 0.02416350844902828
$wsNonLinear.Range("B7").Value = -0.003983050600543842
$wsNonLinear.Range("B8").Value = 0.0032262226905928
$wsNonLinear.Range("B9").Value = 0.03184268139509805
$wsNonLinear.Range("B10").Value = "[1.0, 0.22484813722671818, 0.05817276929431664, 0.08740528076373216, 0.04835190583055679, 0.056385165802227556, 0.20032346151317582, 0.3602105209475373, 0.2110319963757661, 0.045204294879432255, 0.038924525562073646, 0.0446296725009264, 0.06450753686473056, 0.20188637645734123, 0.36354280811727263, 0.20599681409908385, 0.020768151360846195, 0.04526364784172385, 0.0386759456051214, 0.039094157830233706]"
